$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before the existing row 5, shifting rows 5-13 down to 7-15.
$ws.Rows("5:6").Insert()

# New row 5: "Usar lista de contactos no easyphone" (priority 3, status TBD)
$ws.Range("A5").Value = "Usar lista de contactos no easyphone"
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = "TBD"

# New row 6: "Se não existir auricular ligar altifalante" (priority 1, status TBD)
$ws.Range("A6").Value = "Se não existir auricular ligar altifalante"
$ws.Range("B6").Value = 1
$ws.Range("D6").Value = "TBD"

# Append a new row 16 at the bottom: "Ligar altifalante" (priority 1, owner Hugo, status Done)
$ws.Range("A16").Value = "Ligar altifalante"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Hugo"
$ws.Range("D16").Value = "Done"

# Refresh the remembered sort range/state (data is already in the correct
# descending-by-Status order) and move the active selection to D5.
$ws.Range("A2:G16").Sort($ws.Range("D1"), 2, $null, $null, 1)
$ws.Range("D5").Select()
